# Apply hybrid bold + color highlighting to quantitative impact metrics
# (percentages, dollar amounts, large numbers) across several bullet
# points in the document, matching the target OOXML diff.

$d = $word.ActiveDocument

# wdColor-style value for RGB(0x2C, 0x3E, 0x50) stored as a BGR-packed
# long, which is how Word's OLE_COLOR / Font.Color expects values.
$metricColor = 0x2C + (0x3E * 256) + (0x50 * 65536)

function Set-MetricHighlights {
    param(
        [int]$ParaIndex,
        [string[]]$Metrics
    )

    $para = $d.Paragraphs.Item($ParaIndex)
    $paraEnd = $para.Range.End
    $r = $d.Range($para.Range.Start, $paraEnd)

    foreach ($metric in $Metrics) {
        $found = $r.Find.Execute($metric, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $found) {
            throw "Could not find metric '$metric' in paragraph $ParaIndex"
        }
        $r.Font.Bold = 1
        $r.Font.Color = $metricColor
        $r.Collapse(0)
        $r.End = $paraEnd
    }
}

# 1) "• Discovered systematic race coding errors ... from 23% to 64%"
Set-MetricHighlights 9 @("23%", "64%")

# 2) "• Achieved 87% prediction accuracy ... 71%, reducing polling error
#     margins from ±4.2% to ±2.1%"
Set-MetricHighlights 11 @("87%", "71%", "±4.2%", "±2.1%")

# 3) "• Wrote RFP and analyzed bids from 1,200 vendors ..."
Set-MetricHighlights 31 @("1,200")

# 4) "• Created comprehensive meta-analysis framework ... $400M ...
#     now valued at $1B+"
Set-MetricHighlights 46 @('$400M', '$1B')

# 5) "• Algorithm reduced mapping costs by 73.5%, saving campaigns and
#     organizations $4.7M"
Set-MetricHighlights 63 @('73.5%', '$4.7M')

# 6) "• Achieved 87% prediction accuracy for voter turnout vs. industry
#     standard of 71%" (short variant, no polling-error clause)
Set-MetricHighlights 65 @("87%", "71%")
